# Lab01_ReviewReport.xlsx edit: "Tool based analysis cu Sonarlint"
# Fills in the Tool-based Code Analysis sheet with the SonarLint findings,
# fixes two typos/omissions in the Coding Phase Defects sheet, and moves
# the active selection to the Tool-based Code Analysis sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Coding Phase Defects: tidy up two comments, and update the cursor.
# ---------------------------------------------------------------------
$wsCoding = $wb.Worksheets.Item("Coding Phase Defects")

$wsCoding.Range("E10").Value = "In cod nu se verifica daca o comanda este selectata inainte de a apasa pe butonul Ready asa cum ar fi logic,ceea ce duce la o eroare deoarece selectedOrder ar putea fi null."
$wsCoding.Range("E11").Value = "Unii parametri nu sunt denumiti corespunzator ,de exemplu cand vezi variabila ready in cod nu sti daca e un buton pana nu faci actiuni specifice unui buton.Parametrii ar trebui sa aibe nume cat mai sugestive legat de ceea ce reprezinta."

# Row grew taller to fit the longer comment.
$wsCoding.Rows.Item(11).RowHeight = 72

$wsCoding.Range("E12").Select()

# ---------------------------------------------------------------------
# 2. Tool-basedCodeAnalysis: fill in the SonarLint tool-based analysis.
# ---------------------------------------------------------------------
$wsTool = $wb.Worksheets.Item("Tool-basedCodeAnalysis")

$wsTool.Range("D4").Value = "SonarLint"

# Row 10 - Crt. No. 1
$wsTool.Range("C10").Value = "OrdersGUIController.java line 65"
$wsTool.Range("D10").Value = "Metodele are trebui sa nu fie goale"
$wsTool.Range("E10").Value = "OrdersGUIController avea o metoda constructor in care nu se facea nimic"
$wsTool.Range("F10").Value = "A fost adaugat un comentariu inauntrul metodei "
$wsTool.Rows.Item(10).RowHeight = 43.2

# Row 11 - Crt. No. 2
$wsTool.Range("C11").Value = "OrdersGUIController.java line 53,109"
$wsTool.Range("D11").Value = "Membri statici ar trebui sa fie accesati static"
$wsTool.Range("E11").Value = "In metoda setTotalAmount aveam this.totalAmount=totalAmount"
$wsTool.Range("F11").Value = "this a fost inlocuit cu numele clasei"
$wsTool.Rows.Item(11).RowHeight = 57.6

# Row 12 - Crt. No. 3
$wsTool.Range("C12").Value = "OrdersGUIController.java line 61"
$wsTool.Range("D12").Value = "Sectiuni de cod nu ar trebui sa fie comentate"
$wsTool.Range("E12").Value = 'Bucata de cod " // = FXCollections.observableArrayList();   " era comentata'
$wsTool.Range("F12").Value = "Codul comentat a fost sters "
$wsTool.Rows.Item(12).RowHeight = 43.2

# Row 13 - Crt. No. 4
$wsTool.Range("C13").Value = "KitchenGUIController.java line 24"
$wsTool.Range("D13").Value = "Nu ar trebui folositi constructori pentrua instantia String,BigInteger etc."
$wsTool.Range("E13").Value = "extractedTableNumberString era instantiat cu =new String();"
$wsTool.Range("F13").Value = 'new String() a fost inlocuit cu "" reprezentat un string gol'
$wsTool.Rows.Item(13).RowHeight = 72

# Row 14 - Crt. No. 5
$wsTool.Range("C14").Value = "KitchenGUIController.java line 21"
$wsTool.Range("D14").Value = 'Fielduri-le "public static" ar trebui sa fie constante'
$wsTool.Range("E14").Value = "order era doar public si static"
$wsTool.Range("F14").Value = "order a fost facut si final"
$wsTool.Rows.Item(14).RowHeight = 43.2

# Make this the active sheet/tab, scrolled down, with F14 selected.
$wsTool.Activate()
$wsTool.Range("F14").Select()
